# Update player parameter values so the player can't go out of the map.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Block starting at row 25 (Position_X row 27, Position_Y row 28, Size_X row 29)
$ws.Range("B27").Value = 810
$ws.Range("B28").Value = 30
$ws.Range("B29").Value = 300

# Block starting at row 31 (Position_X row 33, Position_Y row 34, Size_X row 35)
$ws.Range("B33").Value = 810
$ws.Range("B34").Value = 30
$ws.Range("B35").Value = 300

# Block starting at row 37 (Position_X row 39, Position_Y row 40, Size_X row 41)
$ws.Range("B39").Value = 810
$ws.Range("B40").Value = 1050
$ws.Range("B41").Value = 300

# Block starting at row 43 (Position_X row 45, Position_Y row 46, Size_X row 47)
$ws.Range("B45").Value = 810
$ws.Range("B46").Value = 1050
$ws.Range("B47").Value = 300

# Scroll the view so row 25 is at the top and select B27, matching the saved view state
$excel.ActiveWindow.ScrollRow = 25
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("B27").Select()
